$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '79.119.68'
$ws.Range('E2').Value = '  +3.12%  '
$ws.Range('D3').Value = '3.174.24'
$ws.Range('E3').Value = '  +3.96%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '629.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.227'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +10.92%  '
$ws.Range('B9').Value = 'XRP'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.584'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.47%  '
$ws.Range('D10').Value = '3.175.89'
$ws.Range('E10').Value = '  +4.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.573'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +31.63%  '
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.41'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.41%  '
$ws.Range('D14').Value = '3.756.62'
$ws.Range('E14').Value = '  +4.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000224'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +17.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '31.71'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.20%  '
$ws.Range('D17').Value = '78.955.21'
$ws.Range('E17').Value = '  +3.17%  '
$ws.Range('D18').Value = '3.179.16'
$ws.Range('E18').Value = '  +5.01%  '
$ws.Range('E19').Value = '  +5.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '427.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +13.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +23.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +13.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.77'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.70%  '
$ws.Range('D25').Value = '3.342.20'
$ws.Range('E25').Value = '  +5.96%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.75'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.49%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '76.62'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0000114'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.88'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.68%  '
$ws.Range('E33').Value = '  +4.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '512.22'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.96'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.131'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +24.72%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '22.82'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.28%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.134'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +19.01%  '
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.398'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '163.63'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '19.97'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '192.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.07%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.39'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.813'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +14.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '42.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.623'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.41%  '
